$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name swaps in shared text (column A) due to re-sort ---
$ws.Range("A20").Value = "Italia"
$ws.Range("A21").Value = "Turquia"
$ws.Range("A66").Value = "Moldavia"
$ws.Range("A67").Value = "Nepal"
$ws.Range("A159").Value = "Principado de Andorra"
$ws.Range("A160").Value = "Guyana"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 18:05"

# --- Update numeric data values (new Covid-19 figures) ---
$ws.Range("B4").Value = 5965815
$ws.Range("C4").Value = 10087
$ws.Range("D4").Value = 3256201
$ws.Range("E4").Value = 2526833
$ws.Range("G4").Value = 377
$ws.Range("H4").Value = 182781
$ws.Range("B6").Value = 3280962
$ws.Range("C6").Value = 49208
$ws.Range("D6").Value = 2501882
$ws.Range("E6").Value = 718917
$ws.Range("G6").Value = 551
$ws.Range("H6").Value = 60163
$ws.Range("B20").Value = 262540
$ws.Range("C20").Value = 1366
$ws.Range("D20").Value = 206329
$ws.Range("E20").Value = 20753
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 35458
$ws.Range("B21").Value = 261194
$ws.Range("D21").Value = 238795
$ws.Range("E21").Value = 16236
$ws.Range("H21").Value = 6163
$ws.Range("B23").Value = 238178
$ws.Range("C23").Value = 606
$ws.Range("E23").Value = 19229
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 9349
$ws.Range("B27").Value = 126057
$ws.Range("C27").Value = 88
$ws.Range("D27").Value = 112165
$ws.Range("E27").Value = 4800
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 9092
$ws.Range("D49").Value = 54971
$ws.Range("E49").Value = 1497
$ws.Range("B62").Value = 39964
$ws.Range("C62").Value = 323
$ws.Range("D62").Value = 36402
$ws.Range("E62").Value = 3271
$ws.Range("G62").Value = 7
$ws.Range("H62").Value = 291
$ws.Range("B64").Value = 35707
$ws.Range("C64").Value = 148
$ws.Range("D64").Value = 33281
$ws.Range("E64").Value = 1904
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 522
$ws.Range("B66").Value = 34982
$ws.Range("C66").Value = 624
$ws.Range("D66").Value = 24156
$ws.Range("E66").Value = 9859
$ws.Range("G66").Value = 7
$ws.Range("H66").Value = 967
$ws.Range("B67").Value = 34418
$ws.Range("C67").Value = 885
$ws.Range("D67").Value = 19504
$ws.Range("E67").Value = 14739
$ws.Range("G67").Value = 11
$ws.Range("H67").Value = 175
$ws.Range("B68").Value = 33016
$ws.Range("C68").Value = 213
$ws.Range("D68").Value = 19296
$ws.Range("E68").Value = 13156
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 564
$ws.Range("D108").Value = 5297
$ws.Range("E108").Value = 26
$ws.Range("B127").Value = 2984
$ws.Range("C127").Value = 13
$ws.Range("E127").Value = 153
$ws.Range("B142").Value = 1813
$ws.Range("C142").Value = 15
$ws.Range("E142").Value = 1131
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 48
$ws.Range("B150").Value = 1484
$ws.Range("C150").Value = 10
$ws.Range("E150").Value = 529
$ws.Range("D152").Value = 192
$ws.Range("E152").Value = 1177
$ws.Range("B159").Value = 1098
$ws.Range("C159").Value = 38
$ws.Range("D159").Value = 893
$ws.Range("E159").Value = 152
$ws.Range("H159").Value = 53
$ws.Range("D160").Value = 523
$ws.Range("E160").Value = 506
$ws.Range("H160").Value = 31
$ws.Range("D169").Value = 345
$ws.Range("E169").Value = 229
